$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.227.22"

$ws.Range("D3").Value = "1.583.75"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'209.80"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "1.805.74"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "1.572.57"
$ws.Range("E13").Value = "  -1.90%  "

$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "'0.518"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").Value = "'64.67"
$ws.Range("E16").Value = "  -0.49%  "

$ws.Range("D17").Value = "26.226.10"
$ws.Range("E17").Value = "  -1.71%  "

$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("D21").Value = "'207.03"
$ws.Range("E21").Value = "  -1.49%  "

$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "'2.22"
$ws.Range("E23").Value = "  -3.43%  "

$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").Value = "'144.72"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("E30").Value = "  -1.69%  "

$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.27"
$ws.Range("E34").Value = "  +6.95%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.289.90"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'0.607"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'5.55"
$ws.Range("E41").Value = "  +2.65%  "

$ws.Range("E42").Value = "  -1.41%  "

$ws.Range("D43").Value = "'2.14"
$ws.Range("E43").Value = "  -2.65%  "

$ws.Range("D44").Value = "'62.36"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("D45").Value = "1.718.61"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").Value = "'88.78"
$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("E47").Value = "  -0.26%  "

$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0987"
$ws.Range("E50").Value = "  -6.79%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.01%  "
